$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 676930.0600000001
$ws.Range("I15").Value = 676930.0600000001
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 2030790.18
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -2030621.18

$ws.Range("H32").Value = 2156.5715
$ws.Range("I32").Value = 2325
$ws.Range("J32").Value = 2089.2
$ws.Range("K32").Value = 2325
$ws.Range("L32").Value = 2089.2
$ws.Range("M32").Value = -1999
$ws.Range("N32").Value = -2741.2

$ws.Range("H62").Value = 3321
$ws.Range("I62").Value = 3046.7
$ws.Range("J62").Value = 3663.875
$ws.Range("K62").Value = 3046.7
$ws.Range("L62").Value = 3663.875
$ws.Range("M62").Value = -2422.7
$ws.Range("N62").Value = -4911.875

$ws.Range("H65").Value = 3321
$ws.Range("I65").Value = 3046.7
$ws.Range("J65").Value = 3663.875
$ws.Range("K65").Value = 15233.5
$ws.Range("L65").Value = 18319.375
$ws.Range("M65").Value = -12113.5
$ws.Range("N65").Value = -24559.375

$ws.Range("H76").Value = 3145.7568
$ws.Range("I76").Value = 3072.6667
$ws.Range("J76").Value = 3748.75
$ws.Range("K76").Value = 3072.6667
$ws.Range("L76").Value = 3748.75
$ws.Range("M76").Value = -2757.6667
$ws.Range("N76").Value = -4378.75

$ws.Range("H79").Value = 3145.7568
$ws.Range("I79").Value = 3072.6667
$ws.Range("J79").Value = 3748.75
$ws.Range("K79").Value = 3072.6667
$ws.Range("L79").Value = 3748.75
$ws.Range("M79").Value = -1980.6667
$ws.Range("N79").Value = -5932.75

$ws.Range("H115").Value = 10000958
$ws.Range("I115").Value = 12500672
$ws.Range("J115").Value = 2102
$ws.Range("K115").Value = 37502016
$ws.Range("L115").Value = 6306
$ws.Range("M115").Value = -37500449
$ws.Range("N115").Value = -9440

$ws.Range("H137").Value = 1638.04
$ws.Range("I137").Value = 835.41174
$ws.Range("J137").Value = 3343.625
$ws.Range("K137").Value = 2506.23522
$ws.Range("L137").Value = 10030.875
$ws.Range("M137").Value = 43.76477999999997
$ws.Range("N137").Value = -15130.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 3830.1
$ws.Range("I26").Value = 2028
$ws.Range("J26").Value = 8035
$ws.Range("K26").Value = 2028
$ws.Range("L26").Value = 8035
$ws.Range("M26").Value = -1698
$ws.Range("N26").Value = -8695

$ws.Range("H74").Value = 1262
$ws.Range("I74").Value = 1300.8462
$ws.Range("J74").Value = 1189.8572
$ws.Range("K74").Value = 1300.8462
$ws.Range("L74").Value = 1189.8572
$ws.Range("M74").Value = -426.8462
$ws.Range("N74").Value = -2937.8572

$ws.Range("H77").Value = 1262
$ws.Range("I77").Value = 1300.8462
$ws.Range("J77").Value = 1189.8572
$ws.Range("K77").Value = 6504.231
$ws.Range("L77").Value = 5949.286
$ws.Range("M77").Value = -2136.231
$ws.Range("N77").Value = -14685.286

$ws.Range("H102").Value = 2849.4546
$ws.Range("I102").Value = 2355.75
$ws.Range("J102").Value = 4166
$ws.Range("K102").Value = 2355.75
$ws.Range("L102").Value = 4166
$ws.Range("M102").Value = -733.75
$ws.Range("N102").Value = -7410

$ws.Range("H122").Value = 1728
$ws.Range("I122").Value = 1722.5
$ws.Range("J122").Value = 1742.6666
$ws.Range("K122").Value = 5167.5
$ws.Range("L122").Value = 5227.9998
$ws.Range("M122").Value = -2717.5
$ws.Range("N122").Value = -10127.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 12228.429
$ws.Range("I33").Value = 2649.75
$ws.Range("J33").Value = 25000
$ws.Range("K33").Value = 2649.75
$ws.Range("L33").Value = 25000
$ws.Range("M33").Value = -2313.75
$ws.Range("N33").Value = -25672

$ws.Range("H94").Value = 1292.3334
$ws.Range("I94").Value = 1388.7142
$ws.Range("J94").Value = 955
$ws.Range("K94").Value = 1388.7142
$ws.Range("L94").Value = 955
$ws.Range("M94").Value = -937.7141999999999

$ws.Range("H99").Value = 11426.546
$ws.Range("I99").Value = 18168.166
$ws.Range("J99").Value = 3336.6
$ws.Range("K99").Value = 18168.166
$ws.Range("L99").Value = 3336.6
$ws.Range("M99").Value = -16670.166
$ws.Range("N99").Value = -6332.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 4376078
$ws.Range("I6").Value = 7000025
$ws.Range("J6").Value = 2833.3333
$ws.Range("K6").Value = 7000025
$ws.Range("L6").Value = 2833.3333
$ws.Range("M6").Value = -6999912

$ws.Range("H62").Value = 4870
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 4870
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 4870
$ws.Range("N62").Value = -6118
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 4870
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 4870
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 24350
$ws.Range("N65").Value = -30590
$ws.Range("M65").ClearContents()

$ws.Range("H99").Value = 1985.7826
$ws.Range("I99").Value = 1697.7222
$ws.Range("J99").Value = 3022.8
$ws.Range("K99").Value = 1697.7222
$ws.Range("L99").Value = 3022.8
$ws.Range("M99").Value = -199.7221999999999
$ws.Range("N99").Value = -6018.8

$ws.Range("H122").Value = 1536.625
$ws.Range("I122").Value = 1242.1428
$ws.Range("J122").Value = 1765.6666
$ws.Range("K122").Value = 3726.4284
$ws.Range("L122").Value = 5296.9998
$ws.Range("M122").Value = -1276.4284
$ws.Range("N122").Value = -10196.9998

$ws.Range("H126").Value = 1985.7826
$ws.Range("I126").Value = 1697.7222
$ws.Range("J126").Value = 3022.8
$ws.Range("K126").Value = 5093.1666
$ws.Range("L126").Value = 9068.400000000001
$ws.Range("M126").Value = -2623.1666
$ws.Range("N126").Value = -14008.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1042.28
$ws.Range("I17").Value = 328.42856
$ws.Range("J17").Value = 1319.8889
$ws.Range("K17").Value = 985.28568
$ws.Range("L17").Value = 3959.6667
$ws.Range("M17").Value = -816.28568
$ws.Range("N17").Value = -4297.6667

$ws.Range("H23").Value = 239.22223
$ws.Range("I23").Value = 280.33334
$ws.Range("J23").Value = 218.66667
$ws.Range("K23").Value = 841.0000200000001
$ws.Range("L23").Value = 656.00001
$ws.Range("M23").Value = -606.0000200000001
$ws.Range("N23").Value = -1126.00001

$ws.Range("H86").Value = 856
$ws.Range("I86").Value = 653.8889
$ws.Range("J86").Value = 1219.8
$ws.Range("K86").Value = 1961.6667
$ws.Range("L86").Value = 3659.4
$ws.Range("M86").Value = -775.6667000000002
$ws.Range("N86").Value = -6031.4

$ws.Range("H89").Value = 856
$ws.Range("I89").Value = 653.8889
$ws.Range("J89").Value = 1219.8
$ws.Range("K89").Value = 5885.0001
$ws.Range("L89").Value = 10978.2
$ws.Range("M89").Value = 42.9998999999998
$ws.Range("N89").Value = -22834.2

$ws.Range("H113").Value = 1815499.5
$ws.Range("I113").Value = 5747667.5
$ws.Range("J113").Value = 652.61536
$ws.Range("K113").Value = 17243002.5
$ws.Range("L113").Value = 1957.84608
$ws.Range("M113").Value = -17240832.5
$ws.Range("N113").Value = -6297.84608

$ws.Range("H132").Value = 1657.7142
$ws.Range("I132").Value = 802
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 7218
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -4688

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3024.652
$ws.Range("I80").Value = 2828.2144
$ws.Range("J80").Value = 3330.2222
$ws.Range("K80").Value = 2828.2144
$ws.Range("L80").Value = 3330.2222
$ws.Range("M80").Value = -1830.2144
$ws.Range("N80").Value = -5326.2222

$ws.Range("H83").Value = 3024.652
$ws.Range("I83").Value = 2828.2144
$ws.Range("J83").Value = 3330.2222
$ws.Range("K83").Value = 14141.072
$ws.Range("L83").Value = 16651.111
$ws.Range("M83").Value = -9149.072
$ws.Range("N83").Value = -26635.111

$ws.Range("H126").Value = 4627.1113
$ws.Range("I126").Value = 4528.8887
$ws.Range("J126").Value = 4725.3335
$ws.Range("K126").Value = 13586.6661
$ws.Range("L126").Value = 14176.0005
$ws.Range("M126").Value = -11116.6661
$ws.Range("N126").Value = -19116.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3532.1765
$ws.Range("I7").Value = 3569.7273
$ws.Range("J7").Value = 3463.3333
$ws.Range("K7").Value = 3569.7273
$ws.Range("L7").Value = 3463.3333
$ws.Range("M7").Value = -3457.7273
$ws.Range("N7").Value = -3687.3333

$ws.Range("H126").Value = 3532.1765
$ws.Range("I126").Value = 3569.7273
$ws.Range("J126").Value = 3463.3333
$ws.Range("K126").Value = 10709.1819
$ws.Range("L126").Value = 10389.9999
$ws.Range("M126").Value = -8239.1819
$ws.Range("N126").Value = -15329.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 10012323
$ws.Range("I62").Value = 16684933
$ws.Range("J62").Value = 3408
$ws.Range("K62").Value = 16684933
$ws.Range("L62").Value = 3408
$ws.Range("M62").Value = -16684309
$ws.Range("N62").Value = -4656

$ws.Range("H64").Value = 12750
$ws.Range("I64").Value = 15000
$ws.Range("J64").Value = 12000
$ws.Range("K64").Value = 15000
$ws.Range("L64").Value = 12000
$ws.Range("M64").Value = -14752

$ws.Range("H65").Value = 10012323
$ws.Range("I65").Value = 16684933
$ws.Range("J65").Value = 3408
$ws.Range("K65").Value = 83424665
$ws.Range("L65").Value = 17040
$ws.Range("M65").Value = -83421545
$ws.Range("N65").Value = -23280

$ws.Range("H67").Value = 12750
$ws.Range("I67").Value = 15000
$ws.Range("J67").Value = 12000
$ws.Range("K67").Value = 15000
$ws.Range("L67").Value = 12000
$ws.Range("M67").Value = -14142

$ws.Range("H122").Value = 1833.6897
$ws.Range("I122").Value = 1247.25
$ws.Range("J122").Value = 3136.889
$ws.Range("K122").Value = 3741.75
$ws.Range("L122").Value = 9410.667000000001
$ws.Range("M122").Value = -1291.75
$ws.Range("N122").Value = -14310.667

$ws.Range("H132").Value = 1355.138
$ws.Range("I132").Value = 1009.18604
$ws.Range("J132").Value = 2346.8667
$ws.Range("K132").Value = 3027.55812
$ws.Range("L132").Value = 7040.6001
$ws.Range("M132").Value = -497.5581200000001
$ws.Range("N132").Value = -12100.6001
